$wb = $excel.ActiveWorkbook

# Sheet "y_fitted_on_begin_2016" -> B2:B14
$ws1 = $wb.Worksheets.Item("y_fitted_on_begin_2016")
$values1 = @(
    30.32009812852694,
    31.9276717723328,
    32.0780328183037,
    32.54379460481448,
    32.25856392075065,
    31.90604200663718,
    31.27255309819001,
    31.56531261738984,
    32.22104961281069,
    31.20113533175558,
    30.215870477588,
    30.25228514594359,
    28.26980085393038
)
for ($i = 0; $i -lt $values1.Length; $i++) {
    $row = $i + 2
    $ws1.Cells.Item($row, 2).Value = $values1[$i]
}

# Sheet "y_fitted_on_begin_2021" -> B2:B19
$ws2 = $wb.Worksheets.Item("y_fitted_on_begin_2021")
$values2 = @(
    30.20339941131409,
    31.86634100755996,
    31.97932853014283,
    32.42541722268849,
    32.1473450054343,
    31.80947517816186,
    31.20390373274412,
    31.49107688713051,
    32.12082585899216,
    31.13591524019607,
    30.19549576837826,
    30.24020168257207,
    28.33883154473909,
    27.82165192186463,
    28.32219059128905,
    27.92401038471845,
    26.55365355853328,
    26.40008389658678
)
for ($i = 0; $i -lt $values2.Length; $i++) {
    $row = $i + 2
    $ws2.Cells.Item($row, 2).Value = $values2[$i]
}
